$wb = $excel.ActiveWorkbook

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 5174.2764
$ws.Range("J112").Value = 5471.386
$ws.Range("L112").Value = 16414.158
$ws.Range("N112").Value = -18630.158

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3366.8333
$ws.Range("I116").Value = 3485.7144
$ws.Range("J116").Value = 3200.4
$ws.Range("K116").Value = 3485.7144
$ws.Range("L116").Value = 3200.4
$ws.Range("M116").Value = -43.71439999999984
$ws.Range("N116").Value = -10084.4

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1485.95
$ws.Range("I137").Value = 1683.5385
$ws.Range("J137").Value = 1390.8148
$ws.Range("K137").Value = 5050.6155
$ws.Range("L137").Value = 4172.4444
$ws.Range("M137").Value = -2500.6155
$ws.Range("N137").Value = -9272.4444

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 44009.39
$ws.Range("I63").Value = 185598
$ws.Range("J63").Value = 4679.222
$ws.Range("K63").Value = 185598
$ws.Range("L63").Value = 4679.222
$ws.Range("M63").Value = -184912
$ws.Range("N63").Value = -6051.222

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 44009.39
$ws.Range("I66").Value = 185598
$ws.Range("J66").Value = 4679.222
$ws.Range("K66").Value = 927990
$ws.Range("L66").Value = 23396.11
$ws.Range("M66").Value = -924558
$ws.Range("N66").Value = -30260.11

# ARM row 113
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 44444
$ws.Range("J113").Value = 44444
$ws.Range("L113").Value = 44444
$ws.Range("N113").Value = -53122

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 39205.15
$ws.Range("I122").Value = 64469.562
$ws.Range("J122").Value = 2456.9092
$ws.Range("K122").Value = 193408.686
$ws.Range("L122").Value = 7370.7276
$ws.Range("M122").Value = -190958.686
$ws.Range("N122").Value = -12270.7276

# ARM row 124
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 67401
$ws.Range("J124").Value = 67401
$ws.Range("L124").Value = 67401
$ws.Range("N124").Value = -77221

# ARM row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 43429
$ws.Range("J135").Value = 43429
$ws.Range("L135").Value = 43429
$ws.Range("N135").Value = -53569

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 72790
$ws.Range("J139").Value = 70638.75
$ws.Range("L139").Value = 70638.75
$ws.Range("N139").Value = -80918.75

# BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 15723.75
$ws.Range("J35").Value = 15723.75
$ws.Range("L35").Value = 15723.75
$ws.Range("N35").Value = -16343.75

# BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 37524.715
$ws.Range("J81").Value = 37524.715
$ws.Range("L81").Value = 37524.715
$ws.Range("N81").Value = -39646.715

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 16139.692
$ws.Range("I82").Value = 6678.5
$ws.Range("J82").Value = 17859.908
$ws.Range("K82").Value = 6678.5
$ws.Range("L82").Value = 17859.908
$ws.Range("M82").Value = -6295.5
$ws.Range("N82").Value = -18625.908

# BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 37524.715
$ws.Range("J84").Value = 37524.715
$ws.Range("L84").Value = 112574.145
$ws.Range("N84").Value = -123182.145

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 16139.692
$ws.Range("I85").Value = 6678.5
$ws.Range("J85").Value = 17859.908
$ws.Range("K85").Value = 6678.5
$ws.Range("L85").Value = 17859.908
$ws.Range("M85").Value = -5352.5
$ws.Range("N85").Value = -20511.908

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 805.2
$ws.Range("I94").Value = 737.75
$ws.Range("J94").Value = 925.1111
$ws.Range("K94").Value = 737.75
$ws.Range("L94").Value = 925.1111
$ws.Range("M94").Value = -286.75
$ws.Range("N94").Value = -1827.1111

# BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 40731.332
$ws.Range("J135").Value = 40731.332
$ws.Range("L135").Value = 40731.332
$ws.Range("N135").Value = -50871.332

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7055.245
$ws.Range("I31").Value = 1822.6666
$ws.Range("J31").Value = 8232.575
$ws.Range("K31").Value = 1822.6666
$ws.Range("L31").Value = 8232.575
$ws.Range("M31").Value = -1527.6666
$ws.Range("N31").Value = -8822.575

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7055.245
$ws.Range("I34").Value = 1822.6666
$ws.Range("J34").Value = 8232.575
$ws.Range("K34").Value = 1822.6666
$ws.Range("L34").Value = 8232.575
$ws.Range("M34").Value = -1620.6666
$ws.Range("N34").Value = -8636.575

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2040.3667
$ws.Range("I99").Value = 1744.5714
$ws.Range("J99").Value = 2130.3914
$ws.Range("K99").Value = 1744.5714
$ws.Range("L99").Value = 2130.3914
$ws.Range("M99").Value = -246.5714
$ws.Range("N99").Value = -5126.3914

# CRP row 116
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2040.3667
$ws.Range("I126").Value = 1744.5714
$ws.Range("J126").Value = 2130.3914
$ws.Range("K126").Value = 5233.7142
$ws.Range("L126").Value = 6391.174199999999
$ws.Range("M126").Value = -2763.7142
$ws.Range("N126").Value = -11331.1742

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 624.9286
$ws.Range("I5").Value = 423.9091
$ws.Range("K5").Value = 1271.7273
$ws.Range("M5").Value = -1159.7273

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 624.9286
$ws.Range("I135").Value = 423.9091
$ws.Range("K135").Value = 3815.1819
$ws.Range("M135").Value = -1280.1819

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1814.125
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 2002.6
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 2002.6
$ws.Range("M102").Value = 122
$ws.Range("N102").Value = -5246.6

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3236.4722
$ws.Range("I7").Value = 3223.4443
$ws.Range("J7").Value = 3275.5557
$ws.Range("K7").Value = 3223.4443
$ws.Range("L7").Value = 3275.5557
$ws.Range("M7").Value = -3111.4443
$ws.Range("N7").Value = -3499.5557

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2911.4167
$ws.Range("I122").Value = 2415.3076
$ws.Range("J122").Value = 3497.7273
$ws.Range("K122").Value = 7245.9228
$ws.Range("L122").Value = 10493.1819
$ws.Range("M122").Value = -4795.9228
$ws.Range("N122").Value = -15393.1819

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3236.4722
$ws.Range("I126").Value = 3223.4443
$ws.Range("J126").Value = 3275.5557
$ws.Range("K126").Value = 9670.332900000001
$ws.Range("L126").Value = 9826.667099999999
$ws.Range("M126").Value = -7200.332900000001
$ws.Range("N126").Value = -14766.6671

# WVR row 75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 59800
$ws.Range("J75").Value = 59800
$ws.Range("L75").Value = 59800
$ws.Range("N75").Value = -61672

# WVR row 78
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 59800
$ws.Range("J78").Value = 59800
$ws.Range("L78").Value = 179400
$ws.Range("N78").Value = -188760
